$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 2492
$ws.Range("K3").Value = 2396
$ws.Range("K4").Value = 504
$ws.Range("K5").Value = 158
$ws.Range("K6").Value = 2997
$ws.Range("K7").Value = 8547

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K3").Value = 28
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 168
$ws.Range("K3").Value = 167
$ws.Range("K7").Value = 568

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 189

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 94
$ws.Range("K3").Value = 121
$ws.Range("K5").Value = 9
$ws.Range("K6").Value = 92
$ws.Range("K7").Value = 335

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 45
$ws.Range("K3").Value = 44
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 71
$ws.Range("K7").Value = 273

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 60
$ws.Range("K7").Value = 202

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 45
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 60
$ws.Range("K4").Value = 32
$ws.Range("K6").Value = 69
$ws.Range("K8").Value = 568
$ws.Range("K10").Value = 48
$ws.Range("K14").Value = 49
$ws.Range("K15").Value = 85
$ws.Range("K16").Value = 25
$ws.Range("K19").Value = 251
$ws.Range("K20").Value = 192
$ws.Range("K23").Value = 77
$ws.Range("K24").Value = 27
$ws.Range("K25").Value = 35
$ws.Range("K26").Value = 10
$ws.Range("K27").Value = 92
$ws.Range("K29").Value = 434
$ws.Range("K33").Value = 335
$ws.Range("K37").Value = 273
$ws.Range("K41").Value = 74
$ws.Range("K42").Value = 300
$ws.Range("K43").Value = 75
$ws.Range("K44").Value = 81
$ws.Range("K46").Value = 18
$ws.Range("K48").Value = 106
$ws.Range("K50").Value = 56
$ws.Range("K52").Value = 233
$ws.Range("K53").Value = 126
$ws.Range("K54").Value = 161
$ws.Range("K55").Value = 94
$ws.Range("K57").Value = 23
$ws.Range("K59").Value = 14
$ws.Range("K63").Value = 29
$ws.Range("K65").Value = 202
$ws.Range("K66").Value = 30
$ws.Range("K67").Value = 335
$ws.Range("K68").Value = 22
$ws.Range("K71").Value = 24
$ws.Range("K76").Value = 123
$ws.Range("K78").Value = 114
$ws.Range("K79").Value = 224
$ws.Range("K80").Value = 29
$ws.Range("K83").Value = 189
$ws.Range("K84").Value = 62
$ws.Range("K85").Value = 412
$ws.Range("K87").Value = 8
$ws.Range("K88").Value = 100
$ws.Range("K89").Value = 113
$ws.Range("K90").Value = 75
$ws.Range("K93").Value = 39
$ws.Range("K94").Value = 102
$ws.Range("K95").Value = 138
$ws.Range("K97").Value = 73
$ws.Range("K98").Value = 51
$ws.Range("K99").Value = 153
$ws.Range("K101").Value = 8547

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 102
$ws.Range("K7").Value = 335

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K5").Value = 1
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 117
$ws.Range("K6").Value = 140
$ws.Range("K7").Value = 434

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 18
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 67
$ws.Range("K7").Value = 251

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 16
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 69

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K4").Value = 5
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 74

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 92
$ws.Range("K4").Value = 11
$ws.Range("K6").Value = 120
$ws.Range("K7").Value = 300

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 34
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K5").Value = 4
$ws.Range("K7").Value = 94

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K3").Value = 9
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 74
$ws.Range("K3").Value = 80
$ws.Range("K7").Value = 224

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 66
$ws.Range("K3").Value = 54
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 192

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 85

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 10

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 14

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K4").Value = 3
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 113

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K4").Value = 12
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 152
$ws.Range("K7").Value = 412

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 233

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 8

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 25
